$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NATMI LR-pair values updated per Dr Hou advice (E/K ligand & receptor-expressing cell counts 1 -> 3,
# with downstream totals/specificities recomputed).

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 98.912777
$ws.Range("H2").Value = 296.738331
$ws.Range("I2").Value = 0.8120825131376513
$ws.Range("J2").Value = 0.8120825131376513
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 68.65869266666667
$ws.Range("N2").Value = 205.976078
$ws.Range("O2").Value = 0.6475952735309433
$ws.Range("P2").Value = 0.6475952735309431
$ws.Range("Q2").Value = 6791.221956849536
$ws.Range("R2").Value = 61120.99761164582
$ws.Range("S2").Value = 0.5259007972250731
$ws.Range("T2").Value = 0.5259007972250731

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 98.912777
$ws.Range("H3").Value = 296.738331
$ws.Range("I3").Value = 0.8120825131376513
$ws.Range("J3").Value = 0.8120825131376513
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 28.23943666666667
$ws.Range("N3").Value = 84.71831
$ws.Range("O3").Value = 0.2663570336431459
$ws.Range("P3").Value = 0.2663570336431459
$ws.Range("Q3").Value = 2793.241101615623
$ws.Range("R3").Value = 25139.16991454061
$ws.Range("S3").Value = 0.2163038892728159
$ws.Range("T3").Value = 0.2163038892728159

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 98.912777
$ws.Range("H4").Value = 296.738331
$ws.Range("I4").Value = 0.8120825131376513
$ws.Range("J4").Value = 0.8120825131376513
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 9.122861666666667
$ws.Range("N4").Value = 27.368585
$ws.Range("O4").Value = 0.08604769282591093
$ws.Range("P4").Value = 0.08604769282591092
$ws.Range("Q4").Value = 902.3675816368484
$ws.Range("R4").Value = 8121.308234731635
$ws.Range("S4").Value = 0.0698778266397624
$ws.Range("T4").Value = 0.06987782663976239

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 17.04862266666667
$ws.Range("H5").Value = 51.14586800000001
$ws.Range("I5").Value = 0.1399706767982279
$ws.Range("J5").Value = 0.1399706767982279
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 68.65869266666667
$ws.Range("N5").Value = 205.976078
$ws.Range("O5").Value = 0.6475952735309433
$ws.Range("P5").Value = 0.6475952735309431
$ws.Range("Q5").Value = 1170.536144060634
$ws.Range("R5").Value = 10534.8252965457
$ws.Range("S5").Value = 0.09064434872745968
$ws.Range("T5").Value = 0.09064434872745966

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 17.04862266666667
$ws.Range("H6").Value = 51.14586800000001
$ws.Range("I6").Value = 0.1399706767982279
$ws.Range("J6").Value = 0.1399706767982279
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 28.23943666666667
$ws.Range("N6").Value = 84.71831
$ws.Range("O6").Value = 0.2663570336431459
$ws.Range("P6").Value = 0.2663570336431459
$ws.Range("Q6").Value = 481.4435000492312
$ws.Range("R6").Value = 4332.99150044308
$ws.Range("S6").Value = 0.0372821742689995
$ws.Range("T6").Value = 0.0372821742689995

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 17.04862266666667
$ws.Range("H7").Value = 51.14586800000001
$ws.Range("I7").Value = 0.1399706767982279
$ws.Range("J7").Value = 0.1399706767982279
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 9.122861666666667
$ws.Range("N7").Value = 27.368585
$ws.Range("O7").Value = 0.08604769282591093
$ws.Range("P7").Value = 0.08604769282591092
$ws.Range("Q7").Value = 155.5322261951978
$ws.Range("R7").Value = 1399.79003575678
$ws.Range("S7").Value = 0.01204415380176878
$ws.Range("T7").Value = 0.01204415380176877

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 5.839988000000001
$ws.Range("H8").Value = 17.519964
$ws.Range("I8").Value = 0.0479468100641207
$ws.Range("J8").Value = 0.04794681006412069
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 68.65869266666667
$ws.Range("N8").Value = 205.976078
$ws.Range("O8").Value = 0.6475952735309433
$ws.Range("P8").Value = 0.6475952735309431
$ws.Range("Q8").Value = 400.9659412690214
$ws.Range("R8").Value = 3608.693471421192
$ws.Range("S8").Value = 0.03105012757841043
$ws.Range("T8").Value = 0.03105012757841042

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 5.839988000000001
$ws.Range("H9").Value = 17.519964
$ws.Range("I9").Value = 0.0479468100641207
$ws.Range("J9").Value = 0.04794681006412069
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 28.23943666666667
$ws.Range("N9").Value = 84.71831
$ws.Range("O9").Value = 0.2663570336431459
$ws.Range("P9").Value = 0.2663570336431459
$ws.Range("Q9").Value = 164.9179712600934
$ws.Range("R9").Value = 1484.26174134084
$ws.Range("S9").Value = 0.01277097010133052
$ws.Range("T9").Value = 0.01277097010133052

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 5.839988000000001
$ws.Range("H10").Value = 17.519964
$ws.Range("I10").Value = 0.0479468100641207
$ws.Range("J10").Value = 0.04794681006412069
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 9.122861666666667
$ws.Range("N10").Value = 27.368585
$ws.Range("O10").Value = 0.08604769282591093
$ws.Range("P10").Value = 0.08604769282591092
$ws.Range("Q10").Value = 53.27740265899335
$ws.Range("R10").Value = 479.49662393094
$ws.Range("S10").Value = 0.004125712384379753
$ws.Range("T10").Value = 0.004125712384379751
